$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged cell A1:C2 holds the sheet's only text; update its content.
$ws.Range("A1").Value = "I am a cell2!"
